# Updates cryptos list values (Price column D, Volume(1h) column E)
# as produced by the scheduled GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    # Force text storage so values like "68.528.54" or padded
    # percentages keep their original string formatting instead
    # of being auto-converted to numbers by Excel.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" "68.528.54"
Set-TextCell "D3" "3.763.97"
Set-TextCell "D5" "594.04"
Set-TextCell "E5" "  -0.57%  "
Set-TextCell "D6" "167.31"
Set-TextCell "E6" "  -1.61%  "
Set-TextCell "D7" "3.763.63"
Set-TextCell "E7" "  -0.53%  "
Set-TextCell "E9" "  -1.06%  "
Set-TextCell "E10" "  -2.67%  "
Set-TextCell "E11" "  -1.48%  "
Set-TextCell "E12" "  -1.04%  "
Set-TextCell "D13" "0.0000261"
Set-TextCell "E13" "  -6.59%  "
Set-TextCell "D14" "36.17"
Set-TextCell "E14" "  -1.57%  "
Set-TextCell "D15" "4.395.28"
Set-TextCell "E15" "  -0.41%  "
Set-TextCell "D16" "3.753.36"
Set-TextCell "E16" "  -0.95%  "
Set-TextCell "D17" "68.503.22"
Set-TextCell "E17" "  +0.82%  "
Set-TextCell "D18" "17.96"
Set-TextCell "E18" "  -4.52%  "
Set-TextCell "E19" "  +0.86%  "
Set-TextCell "D20" "7.00"
Set-TextCell "E20" "  -3.02%  "
Set-TextCell "E21" "  +1.45%  "
Set-TextCell "D22" "465.84"
Set-TextCell "E22" "  -0.50%  "
Set-TextCell "D23" "0.698"
Set-TextCell "E23" "  -3.24%  "
Set-TextCell "D24" "84.12"
Set-TextCell "E24" "  +0.42%  "
Set-TextCell "E25" "  -2.16%  "
Set-TextCell "E26" "  -2.71%  "
Set-TextCell "D27" "11.97"
Set-TextCell "E27" "  -1.64%  "
Set-TextCell "D28" "10.06"
Set-TextCell "E28" "  -4.25%  "
Set-TextCell "E29" "  -0.10%  "
Set-TextCell "D30" "3.910.58"
Set-TextCell "E30" "  -0.48%  "
Set-TextCell "E31" "  -4.78%  "
Set-TextCell "D32" "7.35"
Set-TextCell "E32" "  -3.50%  "
Set-TextCell "D33" "30.04"
Set-TextCell "E33" "  -1.75%  "
Set-TextCell "D34" "2.17"
Set-TextCell "E34" "  -3.13%  "
Set-TextCell "D35" "9.23"
Set-TextCell "E35" "  -0.34%  "
Set-TextCell "D37" "3.715.61"
Set-TextCell "E38" "  -3.56%  "
Set-TextCell "D39" "3.41"
Set-TextCell "E39" "  -8.83%  "
Set-TextCell "E40" "  -0.75%  "
Set-TextCell "E41" "  -0.41%  "
Set-TextCell "D42" "5.78"
Set-TextCell "E42" "  -1.22%  "
Set-TextCell "E44" "  -0.01%  "
Set-TextCell "D45" "44.16"
Set-TextCell "E45" "  +9.21%  "
Set-TextCell "E46" "  -3.54%  "
Set-TextCell "D47" "46.86"
Set-TextCell "E47" "  +2.68%  "
Set-TextCell "D48" "1.92"
Set-TextCell "E48" "  -1.91%  "
Set-TextCell "D49" "8.50"
Set-TextCell "E49" "  -2.28%  "
Set-TextCell "D50" "145.29"
Set-TextCell "E50" "  +2.34%  "
Set-TextCell "D51" "390.32"
Set-TextCell "E51" "  -3.13%  "
